$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.413372394153043
$ws.Range("C2").Value = 0.2014560155030836
$ws.Range("D2").Value = 0.1185488847114016
$ws.Range("E2").Value = 0.1291941946857791
$ws.Range("F2").Value = 1.803637194831531
$ws.Range("J2").Value = 0.1737659992641696
$ws.Range("L2").Value = 0.2725791693644055
$ws.Range("M2").Value = 0.3230976196827768
$ws.Range("O2").Value = 4.72691074783711

$ws.Range("B3").Value = 1.314378415875581
$ws.Range("C3").Value = 0.1887386274948994
$ws.Range("D3").Value = 0.1182101425572206
$ws.Range("E3").Value = 0.1301246217882173
$ws.Range("F3").Value = 1.816263475189196
$ws.Range("J3").Value = 0.1750953665012949
$ws.Range("L3").Value = 0.2687900340464466
$ws.Range("M3").Value = 0.3071849845766295
$ws.Range("O3").Value = 4.765067934788021

$ws.Range("B4").Value = 1.253802580273543
$ws.Range("C4").Value = 0.1808928524788627
$ws.Range("D4").Value = 0.1180283284262913
$ws.Range("E4").Value = 0.1307297229449418
$ws.Range("F4").Value = 1.82502869416156
$ws.Range("J4").Value = 0.1759581931710317
$ws.Range("L4").Value = 0.2665419847510435
$ws.Range("M4").Value = 0.2974856371997419
$ws.Range("O4").Value = 4.791342454177823

$ws.Range("B5").Value = 1.229171044146767
$ws.Range("C5").Value = 0.1776864758290913
$ws.Range("D5").Value = 0.11796084946306
$ws.Range("E5").Value = 0.1309848238460352
$ws.Range("F5").Value = 1.82885520361792
$ws.Range("J5").Value = 0.1763215372721403
$ws.Range("L5").Value = 0.2656457351206072
$ws.Range("M5").Value = 0.2935512272236451
$ws.Range("O5").Value = 4.802764799976615

$ws.Range("B6").Value = 1.225084281259626
$ws.Range("C6").Value = 0.1771535123344705
$ws.Range("D6").Value = 0.1179500448484134
$ws.Range("E6").Value = 0.1310276979949418
$ws.Range("F6").Value = 1.829505969844789
$ws.Range("J6").Value = 0.1763825796479548
$ws.Range("L6").Value = 0.2654981156287874
$ws.Range("M6").Value = 0.2928990261310389
$ws.Range("O6").Value = 4.804704661930828

$ws.Range("B7").Value = 1.253470171585775
$ws.Range("C7").Value = 0.1808496469739822
$ws.Range("D7").Value = 0.1180273915696048
$ws.Range("E7").Value = 0.130733128819851
$ws.Range("F7").Value = 1.825079268927425
$ws.Range("J7").Value = 0.175963045815049
$ws.Range("L7").Value = 0.2665298171014214
$ws.Range("M7").Value = 0.2974325024971378
$ws.Range("O7").Value = 4.791493604314979

$ws.Range("B8").Value = 1.379197253359166
$ws.Range("C8").Value = 0.197078910243107
$ws.Range("D8").Value = 0.1184266714328928
$ws.Range("E8").Value = 0.1295079962000898
$ws.Range("F8").Value = 1.80778054949274
$ws.Range("J8").Value = 0.1742147073576037
$ws.Range("L8").Value = 0.2712564483201874
$ws.Range("M8").Value = 0.3175963537850492
$ws.Range("O8").Value = 4.739476478335888

$ws.Range("B9").Value = 1.627332369568535
$ws.Range("C9").Value = 0.2286013970696388
$ws.Range("D9").Value = 0.1194161389843558
$ws.Range("E9").Value = 0.1273731662755377
$ws.Range("F9").Value = 1.781894509591794
$ws.Range("J9").Value = 0.1711549734841622
$ws.Range("L9").Value = 0.2811439992483855
$ws.Range("M9").Value = 0.357691433683371
$ws.Range("O9").Value = 4.660069162567083

$ws.Range("B10").Value = 1.810546480302321
$ws.Range("C10").Value = 0.251568434234656
$ws.Range("D10").Value = 0.1202675146617622
$ws.Range("E10").Value = 0.1259669410027895
$ws.Range("F10").Value = 1.767779362188676
$ws.Range("J10").Value = 0.1691304891479692
$ws.Range("L10").Value = 0.2887808919491732
$ws.Range("M10").Value = 0.3874765644023555
$ws.Range("O10").Value = 4.615532086818888

$ws.Range("B11").Value = 1.894082098309525
$ws.Range("C11").Value = 0.2619734516113681
$ws.Range("D11").Value = 0.1206815454333423
$ws.Range("E11").Value = 0.1253622411721116
$ws.Range("F11").Value = 1.762423678517862
$ws.Range("J11").Value = 0.1682577555639195
$ws.Range("L11").Value = 0.2923350893888284
$ws.Range("M11").Value = 0.4010956942395083
$ws.Range("O11").Value = 4.598274365031557

$ws.Range("B12").Value = 1.925740974398877
$ws.Range("C12").Value = 0.2659072357292303
$ws.Range("D12").Value = 0.1208421452283801
$ws.Range("E12").Value = 0.1251382753191352
$ws.Range("F12").Value = 1.760548882973495
$ws.Range("J12").Value = 0.1679341875448621
$ws.Range("L12").Value = 0.293692399584998
$ws.Range("M12").Value = 0.4062626991946914
$ws.Range("O12").Value = 4.592171528068832

$ws.Range("B13").Value = 1.9189215499635
$ws.Range("C13").Value = 0.2650603116481989
$ws.Range("D13").Value = 0.1208073879651081
$ws.Range("E13").Value = 0.125186287259158
$ws.Range("F13").Value = 1.760945834938923
$ws.Range("J13").Value = 0.1680035663115085
$ws.Range("L13").Value = 0.2933995727195366
$ws.Range("M13").Value = 0.405149463810524
$ws.Range("O13").Value = 4.593466648450402

$ws.Range("B14").Value = 1.896686186464819
$ws.Range("C14").Value = 0.2622972152882141
$ws.Range("D14").Value = 0.1206946817774082
$ws.Range("E14").Value = 0.1253437148097114
$ws.Range("F14").Value = 1.762266365391469
$ws.Range("J14").Value = 0.1682309968864768
$ws.Range("L14").Value = 0.2924465280186013
$ws.Range("M14").Value = 0.4015205933759418
$ws.Range("O14").Value = 4.597763614164222

$ws.Range("B15").Value = 1.883069688813407
$ws.Range("C15").Value = 0.2606039030640943
$ws.Range("D15").Value = 0.120626141986989
$ws.Range("E15").Value = 0.1254407972385168
$ws.Range("F15").Value = 1.763095193309482
$ws.Range("J15").Value = 0.1683712050420088
$ws.Range("L15").Value = 0.291864243623948
$ws.Range("M15").Value = 0.3992990648018093
$ws.Range("O15").Value = 4.60045194246382

$ws.Range("B16").Value = 1.805090924551394
$ws.Range("C16").Value = 0.250887562125456
$ws.Range("D16").Value = 0.1202409922638097
$ws.Range("E16").Value = 0.1260071628044033
$ws.Range("F16").Value = 1.768150814089253
$ws.Range("J16").Value = 0.1691884931710028
$ws.Range("L16").Value = 0.2885502198542582
$ws.Range("M16").Value = 0.3865879004091468
$ws.Range("O16").Value = 4.616720382366793

$ws.Range("B17").Value = 1.757301111757613
$ws.Range("C17").Value = 0.2449157849360688
$ws.Range("D17").Value = 0.1200115421764352
$ws.Range("E17").Value = 0.1263635652952959
$ws.Range("F17").Value = 1.771525202052516
$ws.Range("J17").Value = 0.169702210433357
$ws.Range("L17").Value = 0.286537619463104
$ws.Range("M17").Value = 0.3788076636251247
$ws.Range("O17").Value = 4.627469921995441

$ws.Range("B18").Value = 1.729831713786552
$ws.Range("C18").Value = 0.2414769605222489
$ws.Range("D18").Value = 0.1198820867817645
$ws.Range("E18").Value = 0.1265718535737945
$ws.Range("F18").Value = 1.773566333034509
$ws.Range("J18").Value = 0.1700022261271545
$ws.Range("L18").Value = 0.2853875743669079
$ws.Range("M18").Value = 0.3743392572044684
$ws.Range("O18").Value = 4.633935342196111

$ws.Range("B19").Value = 1.720534196468975
$ws.Range("C19").Value = 0.2403119504863014
$ws.Range("D19").Value = 0.1198386888111429
$ws.Range("E19").Value = 0.1266429426904505
$ws.Range("F19").Value = 1.774274644115749
$ws.Range("J19").Value = 0.1701045863734443
$ws.Range("L19").Value = 0.2849994888654379
$ws.Range("M19").Value = 0.3728274705714512
$ws.Range("O19").Value = 4.63617293735183

$ws.Range("B20").Value = 1.762386564063149
$ws.Range("C20").Value = 0.2455519080318709
$ws.Range("D20").Value = 0.1200357071152212
$ws.Range("E20").Value = 0.1263252847094813
$ws.Range("F20").Value = 1.771155614253203
$ws.Range("J20").Value = 0.1696470546889923
$ws.Range("L20").Value = 0.2867510836796043
$ws.Range("M20").Value = 0.379635204217692
$ws.Range("O20").Value = 4.626296366330365

$ws.Range("B21").Value = 1.903216565994967
$ws.Range("C21").Value = 0.263108978252717
$ws.Range("D21").Value = 0.1207276830020163
$ws.Range("E21").Value = 0.1252973383525777
$ws.Range("F21").Value = 1.76187433276877
$ws.Range("J21").Value = 0.1681640074071673
$ws.Range("L21").Value = 0.2927261514896458
$ws.Range("M21").Value = 0.4025862186551876
$ws.Range("O21").Value = 4.59648975468491

$ws.Range("B22").Value = 1.995406447849746
$ws.Range("C22").Value = 0.2745463104228918
$ws.Range("D22").Value = 0.1212021544665021
$ws.Range("E22").Value = 0.1246547758332732
$ws.Range("F22").Value = 1.756701984709665
$ws.Range("J22").Value = 0.1672350620805645
$ws.Range("L22").Value = 0.2966976751340411
$ws.Range("M22").Value = 0.4176426495161891
$ws.Range("O22").Value = 4.579529365649563

$ws.Range("B23").Value = 1.946189852775206
$ws.Range("C23").Value = 0.2684454678093573
$ws.Range("D23").Value = 0.1209468961146385
$ws.Range("E23").Value = 0.1249950502021977
$ws.Range("F23").Value = 1.759380779759624
$ws.Range("J23").Value = 0.1677271743508797
$ws.Range("L23").Value = 0.2945719524939676
$ws.Range("M23").Value = 0.4096016639101663
$ws.Range("O23").Value = 4.58835068929946

$ws.Range("B24").Value = 1.760087413501708
$ws.Range("C24").Value = 0.245264334116257
$ws.Range("D24").Value = 0.1200247744872343
$ws.Range("E24").Value = 0.1263425808089313
$ws.Range("F24").Value = 1.77132238982297
$ws.Range("J24").Value = 0.1696719760439134
$ws.Range("L24").Value = 0.2866545546156658
$ws.Range("M24").Value = 0.3792610589048806
$ws.Range("O24").Value = 4.626826041974454

$ws.Range("B25").Value = 1.560041525251336
$ws.Range("C25").Value = 0.2201069874045629
$ws.Range("D25").Value = 0.1191265085268682
$ws.Range("E25").Value = 0.127922137170069
$ws.Range("F25").Value = 1.788036576542098
$ws.Range("J25").Value = 0.1719433703702613
$ws.Range("L25").Value = 0.2784034066236813
$ws.Range("M25").Value = 0.3467864406364001
$ws.Range("O25").Value = 4.679128858749749
